# edit.ps1 -- applies the "Almost final- all functionality added" commit to
# flaskshop/static/uploads/language.xlsx (Sheet1 is a Key/value translation
# table: column B = English key, column C = Hebrew value).
#
# Changes implemented:
#   1. Content: two new translation rows appended to the Key/value table
#      (B112/C112 = "special"/"מיוחד", B113/C113 = "View"/"ראה מוצר").
#   2. Cosmetic: the sheet's "plain" (unstyled) cells move from Arial to
#      Calibri (mirrors the workbook-wide Normal-style font swap seen in the
#      diff), while leaving the two specially-styled rows (the JetBrains Mono
#      "obligatory" labels and the Roboto "On Sale" banner cell) untouched.
#   3. Cosmetic: row 67's height increases (15 -> 15.5pt), and columns B/C
#      get very slightly wider -- a side effect Excel itself produces when
#      the Normal font changes, reproduced here as closely as the host
#      allows.
#   4. View: the active selection moves to C116 (the scroll position itself
#      -- "topLeftCell" -- isn't exposed by this host's window object model,
#      so it can't be round-tripped here).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Append the two new Key/value rows -----------------------------------
$ws.Cells.Item(112, 2).Value = "special"
$ws.Cells.Item(112, 3).Value = "מיוחד"
$ws.Cells.Item(113, 2).Value = "View"
$ws.Cells.Item(113, 3).Value = "ראה מוצר"

# --- 2. Normal font Arial -> Calibri, applied to every pre-existing cell that
#        used the default (un-styled) format. Cells carrying an explicit
#        style are intentionally skipped so their look-and-feel is preserved.
$plainCells = @(
    "B1","C1","A2","B2","C2","A3","B3","C3","A4","C4","A5","C5","B6","C6",
    "B7","C7","B8","C8","B9","C9","B10","C10","B11","C11","B12","C12","C13","B14",
    "C14","C15","B16","C16","C17","B18","C18","B19","C19","B20","C20","B21","C21","B22",
    "C22","B23","C23","C24","B25","C25","B26","C26","B27","C27","C28","B29","C29","B30",
    "C30","B31","C31","B32","C32","B33","C33","B34","C34","B35","C35","B36","C36","B37",
    "C37","C38","B39","C39","B40","C40","B41","C41","B42","C42","B43","C43","B44","C44",
    "B45","C45","B46","C46","C47","B48","C48","B49","C49","C50","B51","C51","B52","C52",
    "B53","C53","B54","C54","B55","C55","B56","C56","B57","C57","B58","C58","B59","C59",
    "B60","C60","B61","C61","B62","C62","B63","C63","B64","C64","C65","B66","C66","C67",
    "B68","C68","B69","C69","B70","C70","B71","C71","B72","C72","B73","C73","B74","C74",
    "B75","C75","C76","B77","C77","B78","C78","B79","C79","B80","C80","B81","C81","B82",
    "C82","B83","C83","C84","B85","C85","C86","B87","C87","C88","B89","C89","B90","C90",
    "B91","C91","B92","C92","B93","C93","B94","C94","B95","C95","B96","C96","B97","C97",
    "B98","C98","C99","B100","C100","B101","C101","B102","C102","B103","C103","B104","C104","B105",
    "C105","B106","C106","B107","C107","C108","B109","C109","B110","C110","B111","C111"
)
foreach ($addr in $plainCells) {
    $ws.Range($addr).Font.Name = "Calibri"
}

# --- 3. Row-height / column-width touch-ups ----------------------------------
$ws.Rows.Item(67).RowHeight = 15.5
$ws.Columns.Item(2).ColumnWidth = 20
$ws.Columns.Item(3).ColumnWidth = 18.142857142857142

# --- 4. Selection moves to C116 ----------------------------------------------
$ws.Range("C116").Select() | Out-Null
